$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new rows above the existing data rows (which hold 16/10/2025 and
# 15/10/2025), pushing them down to the bottom of the table, then fill in
# the newly inserted rows with the additional price-history data points
# (04/11/2025 down through 17/10/2025).
$ws.Rows("2:13").Insert()

# Row 2: 04/11/2025
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "04/11/2025"
$ws.Range("B2").Value = "0.953"
$ws.Range("A2:B2").Style = "Normal"
$ws.Range("C2").Value = "SGD"

# Row 3: 03/11/2025
$ws.Range("A3:B3").NumberFormat = "@"
$ws.Range("A3").Value = "03/11/2025"
$ws.Range("B3").Value = "0.963"
$ws.Range("A3:B3").Style = "Normal"
$ws.Range("C3").Value = "SGD"

# Row 4: 31/10/2025
$ws.Range("A4:B4").NumberFormat = "@"
$ws.Range("A4").Value = "31/10/2025"
$ws.Range("B4").Value = "0.959"
$ws.Range("A4:B4").Style = "Normal"
$ws.Range("C4").Value = "SGD"

# Row 5: 30/10/2025
$ws.Range("A5:B5").NumberFormat = "@"
$ws.Range("A5").Value = "30/10/2025"
$ws.Range("B5").Value = "0.958"
$ws.Range("A5:B5").Style = "Normal"
$ws.Range("C5").Value = "SGD"

# Row 6: 29/10/2025
$ws.Range("A6:B6").NumberFormat = "@"
$ws.Range("A6").Value = "29/10/2025"
$ws.Range("B6").Value = "0.962"
$ws.Range("A6:B6").Style = "Normal"
$ws.Range("C6").Value = "SGD"

# Row 7: 28/10/2025
$ws.Range("A7:B7").NumberFormat = "@"
$ws.Range("A7").Value = "28/10/2025"
$ws.Range("B7").Value = "0.961"
$ws.Range("A7:B7").Style = "Normal"
$ws.Range("C7").Value = "SGD"

# Row 8: 27/10/2025
$ws.Range("A8:B8").NumberFormat = "@"
$ws.Range("A8").Value = "27/10/2025"
$ws.Range("B8").Value = "0.963"
$ws.Range("A8:B8").Style = "Normal"
$ws.Range("C8").Value = "SGD"

# Row 9: 24/10/2025
$ws.Range("A9:B9").NumberFormat = "@"
$ws.Range("A9").Value = "24/10/2025"
$ws.Range("B9").Value = "0.954"
$ws.Range("A9:B9").Style = "Normal"
$ws.Range("C9").Value = "SGD"

# Row 10: 23/10/2025
$ws.Range("A10:B10").NumberFormat = "@"
$ws.Range("A10").Value = "23/10/2025"
$ws.Range("B10").Value = "0.949"
$ws.Range("A10:B10").Style = "Normal"
$ws.Range("C10").Value = "SGD"

# Row 11: 22/10/2025
$ws.Range("A11:B11").NumberFormat = "@"
$ws.Range("A11").Value = "22/10/2025"
$ws.Range("B11").Value = "0.946"
$ws.Range("A11:B11").Style = "Normal"
$ws.Range("C11").Value = "SGD"

# Row 12: 21/10/2025
$ws.Range("A12:B12").NumberFormat = "@"
$ws.Range("A12").Value = "21/10/2025"
$ws.Range("B12").Value = "0.947"
$ws.Range("A12:B12").Style = "Normal"
$ws.Range("C12").Value = "SGD"

# Row 13: 17/10/2025
$ws.Range("A13:B13").NumberFormat = "@"
$ws.Range("A13").Value = "17/10/2025"
$ws.Range("B13").Value = "0.949"
$ws.Range("A13:B13").Style = "Normal"
$ws.Range("C13").Value = "SGD"

